$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing historical rows (data corrections) ---

# Row 92 (2022-05-16)
$ws.Range("F92").Value = 20204953.3
$ws.Range("G92").Value = 7
$ws.Range("H92").Value = 471

# Row 93 (2022-05-17)
$ws.Range("F93").Value = 49555690.7
$ws.Range("G93").Value = 400
$ws.Range("H93").Value = 833

# Row 586 (2024-05-27)
$ws.Range("B586").Value = 1030
$ws.Range("D586").Value = 1034
$ws.Range("E586").Value = 953
$ws.Range("F586").Value = 962965672
$ws.Range("G586").Value = 953186
$ws.Range("H586").Value = 3317

# --- Append new rows of quote data (2024-05-29 .. 2024-06-04) ---

$newRows = @(
    @("2024-05-29", 1001.5, 1034,   1035,   995.5, 586738307.5,  583077,  2155),
    @("2024-05-30", 1010,   1002,   1021.5, 997,   1213093212.5, 1198615, 3094),
    @("2024-05-31", 1012,   1009,   1020,   992,   1184970660.5, 1174635, 2916),
    @("2024-06-03", 1017,   1005.5, 1030.5, 1002,  236324736.5,  5,       1455),
    @("2024-06-04", 975,    1014.5, 1020,   970,   748776944.5,  756508,  3325)
)

$r = 588
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).ClearFormats()

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}
